$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 117, pushing the existing rows
# 117:180 down to 118:181 (mirrors the target diff, which adds one new
# data row and shifts all the following rows down by one).
$ws.Rows("117:117").Insert()

# Populate the newly inserted row 117 with its data.
$ws.Range("A117").Value = 5
$ws.Range("B117").Value = "Macroferia Regional de Talca"
$ws.Range("C117").Value = "Maule"
$ws.Range("D117").Value = 44582
$ws.Range("E117").Value = 7
$ws.Range("F117").Value = 100112045
$ws.Range("G117").Value = "Zapallo"
$ws.Range("H117").Value = "Camote"
$ws.Range("I117").Value = "1a nueva(o)"
$ws.Range("J117").Value = 700
$ws.Range("K117").Value = 250
$ws.Range("L117").Value = 300
$ws.Range("M117").Value = 271
$ws.Range("N117").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O117").Value = "Región del Maule"
$ws.Range("P117").Value = 271
$ws.Range("Q117").Value = 1
$ws.Range("R117").Value = "Hortaliza"
